$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.144.30"
$ws.Range("E2").Value = "  +0.40%  "
$ws.Range("D3").Value = "1.882.91"
$ws.Range("E3").Value = "  -0.91%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.47%  "
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5067"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3855"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.69%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09049"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.86%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.128"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.35%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.60"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.83%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.381"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.32%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.85"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.24%  "
$ws.Range("D14").Value = "1.864.47"
$ws.Range("E14").Value = "  -1.74%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.282"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.28%  "
$ws.Range("E16").Value = "  +0.02%  "
$ws.Range("E17").Value = "  -0.53%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "91.50"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.10%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06596"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.29"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.54%  "
$ws.Range("E21").Value = "  -0.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.142"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.94%  "
$ws.Range("D23").Value = "28.158.35"
$ws.Range("E23").Value = "  +0.30%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.50"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.69%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.263"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.563"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.62%  "
$ws.Range("D27").Value = "2.090.74"
$ws.Range("E27").Value = "  -1.36%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.89"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "156.87"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.15%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "127.24"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.49%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1062"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.06%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.064"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.48%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.630"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.28%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.600"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.26%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.600"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.45%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06620"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.08%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02409"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.29%  "
$ws.Range("E38").Value = "  +1.11%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.297"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.59%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.217"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.92%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6444"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.61%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.56"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.51%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.933"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.20%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6057"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.24"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.15%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.667"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.274"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.00%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.242"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.56%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.010"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.49%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "121.45"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.49%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.72"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.90%  "
